$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates per diff (cryptos price/volume refresh).
# Price-looking numeric strings get NumberFormat "@" first so Excel
# stores them as text (matching the original inlineStr cells) instead
# of silently parsing them into floats and losing formatting like
# trailing zeros (e.g. "1.000", "12.80").

$ws.Range("D2").Value = '30.256.10'
$ws.Range("E2").Value = '  -0.39%  '
$ws.Range("D3").Value = '1.857.91'
$ws.Range("E3").Value = '  -0.95%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.63'
$ws.Range("E5").Value = '  -2.37%  '
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4742'
$ws.Range("E7").Value = '  -1.27%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2747'
$ws.Range("E8").Value = '  -2.35%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06435'
$ws.Range("E9").Value = '  -1.10%  '
$ws.Range("D10").Value = '1.904.99'
$ws.Range("E10").Value = '  +1.56%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07435'
$ws.Range("E11").Value = '  -0.58%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.06'
$ws.Range("E12").Value = '  -3.03%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.987'
$ws.Range("E13").Value = '  -1.63%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '85.49'
$ws.Range("E14").Value = '  -3.06%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6336'
$ws.Range("E15").Value = '  -3.96%  '
$ws.Range("D16").Value = '30.242.86'
$ws.Range("E16").Value = '  -0.35%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.80'
$ws.Range("E18").Value = '  -3.65%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007350'
$ws.Range("E19").Value = '  -2.91%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '225.17'
$ws.Range("E20").Value = '  +2.03%  '
$ws.Range("D21").Value = '2.089.92'
$ws.Range("E21").Value = '  -1.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.117'
$ws.Range("E23").Value = '  -3.39%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.035'
$ws.Range("E24").Value = '  -2.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '167.40'
$ws.Range("E25").Value = '  +0.45%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.251'
$ws.Range("E26").Value = '  -0.85%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.84'
$ws.Range("E27").Value = '  -3.61%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.861'
$ws.Range("E28").Value = '  -5.14%  '
$ws.Range("E29").Value = '  +9.62%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.383'
$ws.Range("E30").Value = '  -5.50%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.237'
$ws.Range("E31").Value = '  -1.58%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.905'
$ws.Range("E32").Value = '  -2.87%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04896'
$ws.Range("E33").Value = '  -2.56%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.147'
$ws.Range("E34").Value = '  -4.17%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7295'
$ws.Range("E35").Value = '  -1.54%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.000'
$ws.Range("E36").Value = '  +0.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.683'
$ws.Range("E37").Value = '  -0.99%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01931'
$ws.Range("E38").Value = '  +6.22%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.625'
$ws.Range("E39").Value = '  +0.35%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9035'
$ws.Range("E40").Value = '  -0.23%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.986'
$ws.Range("E41").Value = '  -3.38%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '105.90'
$ws.Range("E42").Value = '  -0.67%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9951'
$ws.Range("E43").Value = '  -0.96%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4104'
$ws.Range("E44").Value = '  -3.73%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.566'
$ws.Range("E45").Value = '  -5.36%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.063'
$ws.Range("E46").Value = '  -4.45%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '61.51'
$ws.Range("E47").Value = '  -4.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1206'
$ws.Range("E48").Value = '  -5.23%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.739'
$ws.Range("E49").Value = '  -1.96%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.399'
$ws.Range("E50").Value = '  -5.01%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05594'
$ws.Range("E51").Value = '  -2.35%  '
